# MAJ taches et avancement
# Rebuilds the sprint-tracking table: swaps the ITERATION-1 task list for the
# new one, moves the RAF-by-sprint mini table from B11:C16 to M6:N11 (next to
# the chart), re-points the chart series at the new range, and repositions
# the chart further right on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column C is a bit wider now (longer task descriptions). ColumnWidth is
#    quantized to whole pixels by Excel's object model, so 46.25 is the
#    closest achievable value to the target stored width of 47.140625.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 46.25

# ---------------------------------------------------------------------------
# 2. Header row (row 4) is unchanged; row 5-10 hold the task list; row 11
#    carries the totals. Clear the old tail of the sheet (old row 8 SUM
#    formulas + the old B11:C16 "Sprint/RAF" mini table) before writing the
#    new layout so nothing stale is left behind.
# ---------------------------------------------------------------------------
$ws.Range("B8:G16").ClearContents()

# Rows 8-10 are brand new data rows (formerly the SUM-formula row + blank
# rows) - they need the same bordered style ("s=2") already used by rows
# 5-7, so copy that formatting down before filling in the values.
$ws.Range("B7:G7").Copy()
$ws.Range("B8:G10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Task table, rows 5-10.
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "#1"
$ws.Range("C5").Value = "Montée en compétence sur Jmonkey & installation"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 100

$ws.Range("B6").Value = "#2"
$ws.Range("C6").Value = "Création des tâches et graphes agiles"
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = ""

$ws.Range("B7").Value = "#3"
$ws.Range("C7").Value = "Ouvrir une fenêtre avec un monde  16x16 blocs"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = ""

$ws.Range("B8").Value = "#4"
$ws.Range("C8").Value = "Vue à la première personne + imposer gravité"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = ""

$ws.Range("B9").Value = "#5"
$ws.Range("C9").Value = "Intégration système de collisions (solide)"
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = ""

$ws.Range("B10").Value = "#6"
$ws.Range("C10").Value = "Poser des blocs et enlever des blocs"
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = ""

# ---------------------------------------------------------------------------
# 4. Totals row, row 11.
# ---------------------------------------------------------------------------
$ws.Range("E11").Formula = "=SUM(E5:E10)"
$ws.Range("F11").Formula = "=SUM(F5:F10)"
$ws.Range("G11").Value = 90

# ---------------------------------------------------------------------------
# 5. RAF-by-sprint mini table, now at M6:N11 (headers + 5 data points) next
#    to the chart instead of B11:C16.
# ---------------------------------------------------------------------------
$ws.Range("M6").Value = "Sprint"
$ws.Range("N6").Value = "RAF"

$ws.Range("M7").Value = 1
$ws.Range("N7").Value = 100

$ws.Range("M8").Value = 2
$ws.Range("N8").Value = 80

$ws.Range("M9").Value = 3
$ws.Range("N9").Value = 60

$ws.Range("M10").Value = 4
$ws.Range("N10").Value = 30

$ws.Range("M11").Value = 5
$ws.Range("N11").Value = 0

# ---------------------------------------------------------------------------
# 6. Re-point the chart series at the new M/N range.
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart()
$series = $chart.SeriesCollection().Item(1)
$series.Formula = '=SERIES(,Feuil1!$M$7:$M$11,Feuil1!$N$7:$N$11,1)'

# ---------------------------------------------------------------------------
# 7. Move the chart further right/up (from col 4/row 9 to col 15/row 4,
#    anchored with the same pixel offsets as the target layout). Compute the
#    pixel position dynamically from live column widths/row heights so it is
#    correct regardless of the column-C resize above.
# ---------------------------------------------------------------------------
function Get-PixelLeft($col0based, $colOffEmu) {
    $total = 0.0
    for ($i = 1; $i -le $col0based; $i++) {
        $total += $ws.Cells.Item(1, $i).Width()
    }
    return $total + ($colOffEmu / 12700.0)
}

function Get-PixelTop($row0based, $rowOffEmu) {
    $total = 0.0
    for ($i = 1; $i -le $row0based; $i++) {
        $total += $ws.Cells.Item($i, 1).Height()
    }
    return $total + ($rowOffEmu / 12700.0)
}

$newLeft = Get-PixelLeft 15 210911
$newTop = Get-PixelTop 4 34698
$newRight = Get-PixelLeft 18 639536
$newBottom = Get-PixelTop 18 110898

$co.Left = $newLeft
$co.Top = $newTop
$co.Width = $newRight - $newLeft
$co.Height = $newBottom - $newTop

# ---------------------------------------------------------------------------
# 8. Sheet view: drop the old scroll position, zoom back to 100%, select E9.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 100
$ws.Range("E9").Select()
